$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value, derived from the target diff.
$updates = @{
    "G2"  = 3.5934
    "G3"  = 3.3008
    "G4"  = 0.1131
    "G5"  = 0.066
    "G6"  = 0.1108
    "G7"  = 0.0441
    "G8"  = 0.1109
    "G9"  = 0.0674
    "G10" = 0.2542
    "G11" = 0.1812
    "G12" = 3.3341
    "G13" = 1.3333
    "G14" = 1.1453
    "G15" = 1.0004
    "G16" = 0.2728
    "G17" = 0.2563
    "G18" = 0.659
    "G19" = 0.646

    "C20" = 3911.0
    "F20" = 314.0
    "G20" = 0.5927
    "G21" = 0.4664
    "G22" = 0.2838
    "G23" = 0.2725

    "C24" = 6935.0
    "F24" = 504.0
    "G24" = 0.7085
    "G25" = 0.699
    "G26" = 1.0056
    "G27" = 0.9665

    "C28" = 5609.0
    "F28" = 390.0
    "G28" = 0.8486
    "G29" = 0.6419
    "G30" = 0.5675
    "G31" = 0.4523

    "C32" = 30865.0
    "F32" = 1478.0
    "G32" = 3.2162
    "G33" = 3.4605
    "G34" = 1.9287
    "G35" = 1.5421
    "G36" = 1.3617
    "G37" = 0.9535

    "C38" = 23579.0
    "F38" = 1386.0
    "G38" = 2.3912
    "G39" = 1.5447

    "C40" = 34782.0
    "F40" = 1498.0
    "G40" = 2.5822
    "G41" = 1.7916
    "G42" = 1.6502
    "G43" = 1.0333
    "G44" = 1.1184
    "G45" = 0.7144

    "C46" = 26292.0
    "F46" = 1432.0
    "G46" = 2.753
    "G47" = 1.5457
    "G48" = 2.3949
    "G49" = 1.6001
    "G50" = 1.1566
    "G51" = 0.7016

    "C52" = 38856.0
    "F52" = 1570.0
    "G52" = 2.553
    "G53" = 1.827

    "C54" = 37138.0
    "F54" = 1547.0
    "G54" = 2.2739
    "G55" = 1.8458

    "C56" = 34510.0
    "F56" = 1514.0
    "G56" = 2.0377
    "G57" = 1.7439
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
